$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2390.625
$ws.Range("J17").Value = 2390.625
$ws.Range("L17").Value = 7171.875
$ws.Range("N17").Value = -7507.875
$ws.Range("H86").Value = 3293492.2
$ws.Range("I86").Value = 3415.3333
$ws.Range("K86").Value = 3415.3333
$ws.Range("M86").Value = -2292.3333
$ws.Range("H89").Value = 3293492.2
$ws.Range("I89").Value = 3415.3333
$ws.Range("K89").Value = 17076.6665
$ws.Range("M89").Value = -11460.6665
$ws.Range("H100").Value = 11577.786
$ws.Range("I100").Value = 1599.5
$ws.Range("K100").Value = 1599.5
$ws.Range("M100").Value = -1058.5
$ws.Range("H134").Value = 53569.855
$ws.Range("J134").Value = 53569.855
$ws.Range("L134").Value = 53569.855
$ws.Range("N134").Value = -63709.855

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 48808.953
$ws.Range("I2").Value = 53651.316
$ws.Range("K2").Value = 53651.316
$ws.Range("M2").Value = -53538.316
$ws.Range("H97").Value = 1440.5
$ws.Range("I97").Value = 1488.3334
$ws.Range("K97").Value = 1488.3334
$ws.Range("M97").Value = -992.3334
$ws.Range("H116").Value = 48808.953
$ws.Range("I116").Value = 53651.316
$ws.Range("K116").Value = 53651.316
$ws.Range("M116").Value = -51357.316
$ws.Range("H132").Value = 1994.6
$ws.Range("I132").Value = 1995.1613
$ws.Range("K132").Value = 5985.4839
$ws.Range("M132").Value = -3455.4839

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 48808.953
$ws.Range("I3").Value = 53651.316
$ws.Range("K3").Value = 53651.316
$ws.Range("M3").Value = -53537.316

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 295.8
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H31").Value = 34038.12
$ws.Range("I31").Value = 2612.25
$ws.Range("J31").Value = 51995.76
$ws.Range("K31").Value = 2612.25
$ws.Range("L31").Value = 51995.76
$ws.Range("M31").Value = -2317.25
$ws.Range("N31").Value = -52585.76
$ws.Range("H34").Value = 34038.12
$ws.Range("I34").Value = 2612.25
$ws.Range("J34").Value = 51995.76
$ws.Range("K34").Value = 2612.25
$ws.Range("L34").Value = 51995.76
$ws.Range("M34").Value = -2410.25
$ws.Range("N34").Value = -52399.76
$ws.Range("H62").Value = 3497.5
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3497.5
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H134").Value = 201661.81
$ws.Range("I134").Value = 1730.6595
$ws.Range("J134").Value = 3333916.8
$ws.Range("K134").Value = 5191.9785
$ws.Range("L134").Value = 10001750.4
$ws.Range("M134").Value = -2656.9785
$ws.Range("N134").Value = -10006820.4

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1509505.4
$ws.Range("I4").Value = 1782951.8
$ws.Range("K4").Value = 5348855.4
$ws.Range("M4").Value = -5348743.4
$ws.Range("H5").Value = 501207.2
$ws.Range("J5").Value = 1112274.6
$ws.Range("L5").Value = 3336823.8
$ws.Range("N5").Value = -3337047.8
$ws.Range("H101").Value = 4675.5
$ws.Range("J101").Value = 3325
$ws.Range("L101").Value = 9975
$ws.Range("N101").Value = -14843
$ws.Range("H131").Value = 3742.238
$ws.Range("I131").Value = 1206.75
$ws.Range("K131").Value = 3620.25
$ws.Range("M131").Value = 1419.75
$ws.Range("H135").Value = 501207.2
$ws.Range("J135").Value = 1112274.6
$ws.Range("L135").Value = 10010471.4
$ws.Range("N135").Value = -10015541.4

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 58654.285
$ws.Range("J32").Value = 58654.285
$ws.Range("L32").Value = 58654.285
$ws.Range("N32").Value = -59246.285
$ws.Range("H113").Value = 2831.3333
$ws.Range("I113").Value = 1515.3077
$ws.Range("K113").Value = 1515.3077
$ws.Range("M113").Value = 654.6922999999999
$ws.Range("H122").Value = 6001
$ws.Range("I122").Value = 5669
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 17007
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -14557
$ws.Range("N122").Value = -23650
$ws.Range("H126").Value = 3040.121
$ws.Range("I126").Value = 2576.2
$ws.Range("K126").Value = 7728.599999999999
$ws.Range("M126").Value = -5258.599999999999
$ws.Range("H132").Value = 35285.25
$ws.Range("I132").Value = 4805.909
$ws.Range("J132").Value = 102339.8
$ws.Range("K132").Value = 14417.727
$ws.Range("L132").Value = 307019.4
$ws.Range("M132").Value = -11887.727
$ws.Range("N132").Value = -312079.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2010744.4
$ws.Range("J7").Value = 2505929.5
$ws.Range("L7").Value = 2505929.5
$ws.Range("N7").Value = -2506153.5
$ws.Range("H22").Value = 1798.6
$ws.Range("I22").Value = 1123.375
$ws.Range("K22").Value = 1123.375
$ws.Range("M22").Value = -828.375
$ws.Range("H27").Value = 1798.6
$ws.Range("I27").Value = 1123.375
$ws.Range("K27").Value = 1123.375
$ws.Range("M27").Value = -1016.375
$ws.Range("H55").Value = 1257.1875
$ws.Range("I55").Value = 313.375
$ws.Range("J55").Value = 2201
$ws.Range("K55").Value = 313.375
$ws.Range("L55").Value = 2201
$ws.Range("M55").Value = -140.375
$ws.Range("N55").Value = -2547
$ws.Range("H61").Value = 4718.533
$ws.Range("I61").Value = 4416.273
$ws.Range("J61").Value = 5549.75
$ws.Range("K61").Value = 4416.273
$ws.Range("L61").Value = 5549.75
$ws.Range("M61").Value = -4214.273
$ws.Range("N61").Value = -5953.75
$ws.Range("H82").Value = 897.6667
$ws.Range("I82").Value = 878.8333
$ws.Range("K82").Value = 878.8333
$ws.Range("M82").Value = -517.8333
$ws.Range("H85").Value = 897.6667
$ws.Range("I85").Value = 878.8333
$ws.Range("K85").Value = 878.8333
$ws.Range("M85").Value = 369.1667
$ws.Range("H93").Value = 2828
$ws.Range("I93").Value = 2322.625
$ws.Range("K93").Value = 2322.625
$ws.Range("M93").Value = -1074.625
$ws.Range("H113").Value = 4718.533
$ws.Range("I113").Value = 4416.273
$ws.Range("J113").Value = 5549.75
$ws.Range("K113").Value = 4416.273
$ws.Range("L113").Value = 5549.75
$ws.Range("M113").Value = -2246.273
$ws.Range("N113").Value = -9889.75
$ws.Range("H122").Value = 1054861
$ws.Range("I122").Value = 590579.9
$ws.Range("K122").Value = 1771739.7
$ws.Range("M122").Value = -1769289.7
$ws.Range("H126").Value = 2010744.4
$ws.Range("J126").Value = 2505929.5
$ws.Range("L126").Value = 7517788.5
$ws.Range("N126").Value = -7522728.5
$ws.Range("H136").Value = 957238.7
$ws.Range("I136").Value = 1057158.4
$ws.Range("K136").Value = 3171475.2
$ws.Range("M136").Value = -3168925.2

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 30802.527
$ws.Range("J132").Value = 95626.27
$ws.Range("L132").Value = 286878.81
$ws.Range("N132").Value = -291938.81
$ws.Range("H136").Value = 9335085
$ws.Range("I136").Value = 10730864
$ws.Range("J136").Value = 402095.8
$ws.Range("K136").Value = 32192592
$ws.Range("L136").Value = 1206287.4
$ws.Range("M136").Value = -32190042
$ws.Range("N136").Value = -1211387.4
